$wb = $excel.ActiveWorkbook

# "Ready for handoff" -> "In Translation" (shared string used by the Status
# cells on all three sheets: Overview!E2/F2, zh-cn!C2, de-de!C2).
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"

# Column widths for the same Status columns shrank (narrower text no longer
# needs as much room). ColumnWidth is expressed in characters and is
# quantized internally, so 12.5 is the closest achievable value to the
# target stored width.
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
